$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclaimer text: date changes from 2021-04-05 to 2021-04-06
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."

# Update Weight (column D) and Percent Change (column E) values for rows 2-15
$ws.Range("D2").Value = 0.05710333091413305
$ws.Range("E2").Value = 0.0003778575477044743

$ws.Range("D3").Value = 0.0238357588774632
$ws.Range("E3").Value = -0.008712871287128721

$ws.Range("D4").Value = 0.03193856833801412
$ws.Range("E4").Value = 0.004750142504275079

$ws.Range("D5").Value = 0.03175921015240153
$ws.Range("E5").Value = -0.003227107704719567

$ws.Range("D6").Value = 0.035253997656274
$ws.Range("E6").Value = -0.007114987376635407

$ws.Range("D7").Value = 0.01904288336219191
$ws.Range("E7").Value = 0.001598691303862809

$ws.Range("D8").Value = 0.004816374134702111
$ws.Range("E8").Value = 0

$ws.Range("D9").Value = 0.006810216822132342
$ws.Range("E9").Value = 0.004752475247524757

$ws.Range("D10").Value = 0.06942375485215306
$ws.Range("E10").Value = 0.009324009324009452

$ws.Range("D11").Value = 0.06958558178654037
$ws.Range("E11").Value = 0.008139534883720989

$ws.Range("D12").Value = 0.1477048372464379
$ws.Range("E12").Value = 0.006792783580454431

$ws.Range("D13").Value = 0.388724479091741
$ws.Range("E13").Value = 0.004865534324133058

$ws.Range("D14").Value = 0.1140010067658156
$ws.Range("E14").Value = -0.003229412495342254

$ws.Range("E15").Value = 0.003415323736299714

# Re-apply sheet protection (the sheet was protected before this edit)
$ws.Protect()
